$d = $word.ActiveDocument

# Create the three new MSC_Join_* paragraph styles, each based on the
# existing MSCJoin ("MSC_Join") style, mirroring the per-language
# MSC_Paragraph_* variants (A = default font, B = Simplified Chinese,
# C = Korean).
$sA = $d.Styles.Add("MSC_Join_A", 1)
$sA.BaseStyle = $d.Styles("MSCJoin")

$sB = $d.Styles.Add("MSC_Join_B", 1)
$sB.BaseStyle = $d.Styles("MSCJoin")
$sB.Font.NameAscii = "Noto Sans CJK SC"
$sB.Font.NameFarEast = "Noto Sans CJK SC"
$sB.Font.NameOther = "Noto Sans CJK SC"
$sB.Font.NameBi = "Noto Sans CJK SC"

$sC = $d.Styles.Add("MSC_Join_C", 1)
$sC.BaseStyle = $d.Styles("MSCJoin")
$sC.Font.NameAscii = "Noto Sans CJK KR"
$sC.Font.NameFarEast = "Noto Sans CJK KR"
$sC.Font.NameOther = "Noto Sans CJK KR"
$sC.Font.NameBi = "Noto Sans CJK KR"

# Re-point every "MSCJoin"-styled paragraph to the matching language
# specific join style, based on which MSC_Paragraph_* style precedes it.
$lastLang = ""
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $sname = $p.Style.NameLocal
    if ($sname -eq "MSC_Paragraph_A") {
        $lastLang = "A"
    } elseif ($sname -eq "MSC_Paragraph_B") {
        $lastLang = "B"
    } elseif ($sname -eq "MSC_Paragraph_C") {
        $lastLang = "C"
    } elseif ($sname -eq "MSC_Join") {
        if ($lastLang -eq "A") {
            $p.Style = "MSC_Join_A"
        } elseif ($lastLang -eq "B") {
            $p.Style = "MSC_Join_B"
        } elseif ($lastLang -eq "C") {
            $p.Style = "MSC_Join_C"
        }
    }
}
